# Weekly BRVM data refresh (GitHub Actions automation) for data/recommandations.xlsx
# - 'Recommandations': updated Jours en Hausse/Baisse, Variation Totale/Derniere Variation,
#   re-ranked rows by Variation Totale (%), dropped FILTISAC CI (FTSC) which fell out of the list
# - 'Top_YTD': refreshed YTD progression figures

$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")

# Row 2: BRVM - SERVICES PUBLICS
$wsReco.Cells.Item(2, 1).Value = 'BRVM - SERVICES PUBLICS'
$wsReco.Cells.Item(2, 2).Value = 0
$wsReco.Cells.Item(2, 3).Value = 8
$wsReco.Cells.Item(2, 4).Value = 3439.1
$wsReco.Cells.Item(2, 5).Value = 114.42
$wsReco.Cells.Item(2, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(2, 7).Value = '➖ Neutre'

# Row 3: SAFCA CI
$wsReco.Cells.Item(3, 1).Value = 'SAFCA CI'
$wsReco.Cells.Item(3, 2).Value = 0
$wsReco.Cells.Item(3, 3).Value = 4
$wsReco.Cells.Item(3, 4).Value = 2745
$wsReco.Cells.Item(3, 5).Value = 680
$wsReco.Cells.Item(3, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(3, 7).Value = '➖ Neutre'

# Row 4: CFAO MOTORS CI
$wsReco.Cells.Item(4, 1).Value = 'CFAO MOTORS CI'
$wsReco.Cells.Item(4, 2).Value = 0
$wsReco.Cells.Item(4, 3).Value = 4
$wsReco.Cells.Item(4, 4).Value = 2705
$wsReco.Cells.Item(4, 5).Value = 680
$wsReco.Cells.Item(4, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(4, 7).Value = '➖ Neutre'

# Row 5: BRVM - AUTRES SECTEURS
$wsReco.Cells.Item(5, 1).Value = 'BRVM - AUTRES SECTEURS'
$wsReco.Cells.Item(5, 2).Value = 0
$wsReco.Cells.Item(5, 3).Value = 4
$wsReco.Cells.Item(5, 4).Value = 2657.01
$wsReco.Cells.Item(5, 5).Value = 663.72
$wsReco.Cells.Item(5, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(5, 7).Value = '➖ Neutre'

# Row 6: UNIWAX CI
$wsReco.Cells.Item(6, 1).Value = 'UNIWAX CI'
$wsReco.Cells.Item(6, 2).Value = 0
$wsReco.Cells.Item(6, 3).Value = 4
$wsReco.Cells.Item(6, 4).Value = 2370
$wsReco.Cells.Item(6, 5).Value = 590
$wsReco.Cells.Item(6, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(6, 7).Value = '➖ Neutre'

# Row 7: NEI-CEDA CI
$wsReco.Cells.Item(7, 1).Value = 'NEI-CEDA CI'
$wsReco.Cells.Item(7, 2).Value = 0
$wsReco.Cells.Item(7, 3).Value = 4
$wsReco.Cells.Item(7, 4).Value = 2365
$wsReco.Cells.Item(7, 5).Value = 595
$wsReco.Cells.Item(7, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(7, 7).Value = '➖ Neutre'

# Row 8: AIR LIQUIDE CI
$wsReco.Cells.Item(8, 1).Value = 'AIR LIQUIDE CI'
$wsReco.Cells.Item(8, 2).Value = 0
$wsReco.Cells.Item(8, 3).Value = 4
$wsReco.Cells.Item(8, 4).Value = 2170
$wsReco.Cells.Item(8, 5).Value = 540
$wsReco.Cells.Item(8, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(8, 7).Value = '➖ Neutre'

# Row 9: SETAO CI
$wsReco.Cells.Item(9, 1).Value = 'SETAO CI'
$wsReco.Cells.Item(9, 2).Value = 0
$wsReco.Cells.Item(9, 3).Value = 4
$wsReco.Cells.Item(9, 4).Value = 2145
$wsReco.Cells.Item(9, 5).Value = 520
$wsReco.Cells.Item(9, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(9, 7).Value = '➖ Neutre'

# Row 10: SUCRIVOIRE
$wsReco.Cells.Item(10, 1).Value = 'SUCRIVOIRE'
$wsReco.Cells.Item(10, 2).Value = 0
$wsReco.Cells.Item(10, 3).Value = 2
$wsReco.Cells.Item(10, 4).Value = 1915
$wsReco.Cells.Item(10, 5).Value = 970
$wsReco.Cells.Item(10, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(10, 7).Value = '➖ Neutre'

# Row 11: BRVM - TRANSPORT
$wsReco.Cells.Item(11, 1).Value = 'BRVM - TRANSPORT'
$wsReco.Cells.Item(11, 2).Value = 0
$wsReco.Cells.Item(11, 3).Value = 4
$wsReco.Cells.Item(11, 4).Value = 1562.5
$wsReco.Cells.Item(11, 5).Value = 378.69
$wsReco.Cells.Item(11, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(11, 7).Value = '➖ Neutre'

# Row 12: BRVM - DISTRIBUTION
$wsReco.Cells.Item(12, 1).Value = 'BRVM - DISTRIBUTION'
$wsReco.Cells.Item(12, 2).Value = 0
$wsReco.Cells.Item(12, 3).Value = 4
$wsReco.Cells.Item(12, 4).Value = 1481.94
$wsReco.Cells.Item(12, 5).Value = 366.28
$wsReco.Cells.Item(12, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(12, 7).Value = '➖ Neutre'

# Row 13: BRVM - AGRICULTURE
$wsReco.Cells.Item(13, 1).Value = 'BRVM - AGRICULTURE'
$wsReco.Cells.Item(13, 2).Value = 0
$wsReco.Cells.Item(13, 3).Value = 4
$wsReco.Cells.Item(13, 4).Value = 1332.98
$wsReco.Cells.Item(13, 5).Value = 329.97
$wsReco.Cells.Item(13, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(13, 7).Value = '➖ Neutre'

# Row 14: BRVM - INDUSTRIE
$wsReco.Cells.Item(14, 1).Value = 'BRVM - INDUSTRIE'
$wsReco.Cells.Item(14, 2).Value = 0
$wsReco.Cells.Item(14, 3).Value = 4
$wsReco.Cells.Item(14, 4).Value = 774.8
$wsReco.Cells.Item(14, 5).Value = 193.92
$wsReco.Cells.Item(14, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(14, 7).Value = '➖ Neutre'

# Row 15: BRVM-PRINCIPAL
$wsReco.Cells.Item(15, 1).Value = 'BRVM-PRINCIPAL'
$wsReco.Cells.Item(15, 2).Value = 0
$wsReco.Cells.Item(15, 3).Value = 4
$wsReco.Cells.Item(15, 4).Value = 714.0700000000001
$wsReco.Cells.Item(15, 5).Value = 177.83
$wsReco.Cells.Item(15, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(15, 7).Value = '➖ Neutre'

# Row 16: BRVM - CONSOMMATION DE BASE
$wsReco.Cells.Item(16, 1).Value = 'BRVM - CONSOMMATION DE BASE'
$wsReco.Cells.Item(16, 2).Value = 0
$wsReco.Cells.Item(16, 3).Value = 4
$wsReco.Cells.Item(16, 4).Value = 680.67
$wsReco.Cells.Item(16, 5).Value = 170.65
$wsReco.Cells.Item(16, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(16, 7).Value = '➖ Neutre'

# Row 17: BRVM - INDUSTRIELS
$wsReco.Cells.Item(17, 1).Value = 'BRVM - INDUSTRIELS'
$wsReco.Cells.Item(17, 2).Value = 0
$wsReco.Cells.Item(17, 3).Value = 4
$wsReco.Cells.Item(17, 4).Value = 587.3099999999999
$wsReco.Cells.Item(17, 5).Value = 141.89
$wsReco.Cells.Item(17, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(17, 7).Value = '➖ Neutre'

# Row 18: BRVM-PRESTIGE
$wsReco.Cells.Item(18, 1).Value = 'BRVM-PRESTIGE'
$wsReco.Cells.Item(18, 2).Value = 0
$wsReco.Cells.Item(18, 3).Value = 4
$wsReco.Cells.Item(18, 4).Value = 524.85
$wsReco.Cells.Item(18, 5).Value = 131.11
$wsReco.Cells.Item(18, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(18, 7).Value = '➖ Neutre'

# Row 19: BRVM - FINANCES
$wsReco.Cells.Item(19, 1).Value = 'BRVM - FINANCES'
$wsReco.Cells.Item(19, 2).Value = 0
$wsReco.Cells.Item(19, 3).Value = 4
$wsReco.Cells.Item(19, 4).Value = 494.79
$wsReco.Cells.Item(19, 5).Value = 123.75
$wsReco.Cells.Item(19, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(19, 7).Value = '➖ Neutre'

# Row 20: BRVM - SERVICES FINANCIERS
$wsReco.Cells.Item(20, 1).Value = 'BRVM - SERVICES FINANCIERS'
$wsReco.Cells.Item(20, 2).Value = 0
$wsReco.Cells.Item(20, 3).Value = 4
$wsReco.Cells.Item(20, 4).Value = 486.27
$wsReco.Cells.Item(20, 5).Value = 121.62
$wsReco.Cells.Item(20, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(20, 7).Value = '➖ Neutre'

# Row 21: BRVM - ENERGIE
$wsReco.Cells.Item(21, 1).Value = 'BRVM - ENERGIE'
$wsReco.Cells.Item(21, 2).Value = 0
$wsReco.Cells.Item(21, 3).Value = 4
$wsReco.Cells.Item(21, 4).Value = 443.07
$wsReco.Cells.Item(21, 5).Value = 108.74
$wsReco.Cells.Item(21, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(21, 7).Value = '➖ Neutre'

# Row 22: BRVM - CONSOMMATION DISCRETIONNAIRE
$wsReco.Cells.Item(22, 1).Value = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$wsReco.Cells.Item(22, 2).Value = 0
$wsReco.Cells.Item(22, 3).Value = 4
$wsReco.Cells.Item(22, 4).Value = 428.86
$wsReco.Cells.Item(22, 5).Value = 106.32
$wsReco.Cells.Item(22, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(22, 7).Value = '➖ Neutre'

# Row 23: BRVM - TELECOMMUNICATIONS
$wsReco.Cells.Item(23, 1).Value = 'BRVM - TELECOMMUNICATIONS'
$wsReco.Cells.Item(23, 2).Value = 0
$wsReco.Cells.Item(23, 3).Value = 4
$wsReco.Cells.Item(23, 4).Value = 389.05
$wsReco.Cells.Item(23, 5).Value = 96.03
$wsReco.Cells.Item(23, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(23, 7).Value = '➖ Neutre'

# Row 24: UNILEVER CI (UNLC)
$wsReco.Cells.Item(24, 1).Value = 'UNILEVER CI (UNLC)'
$wsReco.Cells.Item(24, 2).Value = 3
$wsReco.Cells.Item(24, 3).Value = 0
$wsReco.Cells.Item(24, 4).Value = 18.55
$wsReco.Cells.Item(24, 5).Value = 7.49
$wsReco.Cells.Item(24, 6).Value = '🟢 Achat'
$wsReco.Cells.Item(24, 7).Value = '✅ Renforcer'

# Row 25: NSIA BANQUE COTE D'IVOIRE (NSBC)
$wsReco.Cells.Item(25, 1).Value = 'NSIA BANQUE COTE D''IVOIRE (NSBC)'
$wsReco.Cells.Item(25, 2).Value = 2
$wsReco.Cells.Item(25, 3).Value = 0
$wsReco.Cells.Item(25, 4).Value = 7.66
$wsReco.Cells.Item(25, 5).Value = 3.03
$wsReco.Cells.Item(25, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(25, 7).Value = '➖ Neutre'

# Row 26: ECOBANK COTE D''IVOIRE (ECOC)
$wsReco.Cells.Item(26, 1).Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$wsReco.Cells.Item(26, 2).Value = 2
$wsReco.Cells.Item(26, 3).Value = 0
$wsReco.Cells.Item(26, 4).Value = 7.44
$wsReco.Cells.Item(26, 5).Value = 4.3
$wsReco.Cells.Item(26, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(26, 7).Value = '➖ Neutre'

# Row 27: SETAO CI (STAC)
$wsReco.Cells.Item(27, 1).Value = 'SETAO CI (STAC)'
$wsReco.Cells.Item(27, 2).Value = 2
$wsReco.Cells.Item(27, 3).Value = 1
$wsReco.Cells.Item(27, 4).Value = 6.03
$wsReco.Cells.Item(27, 5).Value = 3.85
$wsReco.Cells.Item(27, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(27, 7).Value = '👀 À surveiller'

# Row 28: CIE CI (CIEC)
$wsReco.Cells.Item(28, 1).Value = 'CIE CI (CIEC)'
$wsReco.Cells.Item(28, 2).Value = 1
$wsReco.Cells.Item(28, 3).Value = 0
$wsReco.Cells.Item(28, 4).Value = 4
$wsReco.Cells.Item(28, 5).Value = 4
$wsReco.Cells.Item(28, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(28, 7).Value = '➖ Neutre'

# Row 29: SUCRIVOIRE (SCRC)
$wsReco.Cells.Item(29, 1).Value = 'SUCRIVOIRE (SCRC)'
$wsReco.Cells.Item(29, 2).Value = 1
$wsReco.Cells.Item(29, 3).Value = 0
$wsReco.Cells.Item(29, 4).Value = 3.09
$wsReco.Cells.Item(29, 5).Value = 3.09
$wsReco.Cells.Item(29, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(29, 7).Value = '➖ Neutre'

# Row 30: BANK OF AFRICA BN (BOAB)
$wsReco.Cells.Item(30, 1).Value = 'BANK OF AFRICA BN (BOAB)'
$wsReco.Cells.Item(30, 2).Value = 1
$wsReco.Cells.Item(30, 3).Value = 0
$wsReco.Cells.Item(30, 4).Value = 2.81
$wsReco.Cells.Item(30, 5).Value = 2.81
$wsReco.Cells.Item(30, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(30, 7).Value = '➖ Neutre'

# Row 31: CFAO MOTORS CI (CFAC)
$wsReco.Cells.Item(31, 1).Value = 'CFAO MOTORS CI (CFAC)'
$wsReco.Cells.Item(31, 2).Value = 1
$wsReco.Cells.Item(31, 3).Value = 1
$wsReco.Cells.Item(31, 4).Value = 0.82
$wsReco.Cells.Item(31, 5).Value = 3.03
$wsReco.Cells.Item(31, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(31, 7).Value = '👀 À surveiller'

# Row 32: BERNABE CI (BNBC)
$wsReco.Cells.Item(32, 1).Value = 'BERNABE CI (BNBC)'
$wsReco.Cells.Item(32, 2).Value = 1
$wsReco.Cells.Item(32, 3).Value = 1
$wsReco.Cells.Item(32, 4).Value = 0.5
$wsReco.Cells.Item(32, 5).Value = 3.02
$wsReco.Cells.Item(32, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(32, 7).Value = '👀 À surveiller'

# Row 33: SAFCA CI (SAFC)
$wsReco.Cells.Item(33, 1).Value = 'SAFCA CI (SAFC)'
$wsReco.Cells.Item(33, 2).Value = 1
$wsReco.Cells.Item(33, 3).Value = 1
$wsReco.Cells.Item(33, 4).Value = 0.2
$wsReco.Cells.Item(33, 5).Value = 4.55
$wsReco.Cells.Item(33, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(33, 7).Value = '👀 À surveiller'

# Row 34: TOTAL
$wsReco.Cells.Item(34, 1).Value = 'TOTAL'
$wsReco.Cells.Item(34, 2).Value = 0
$wsReco.Cells.Item(34, 3).Value = 4
$wsReco.Cells.Item(34, 4).Value = 0
$wsReco.Cells.Item(34, 5).Value = 0
$wsReco.Cells.Item(34, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(34, 7).Value = '➖ Neutre'

# Row 35: ONATEL BF (ONTBF)
$wsReco.Cells.Item(35, 1).Value = 'ONATEL BF (ONTBF)'
$wsReco.Cells.Item(35, 2).Value = 1
$wsReco.Cells.Item(35, 3).Value = 1
$wsReco.Cells.Item(35, 4).Value = -0.02
$wsReco.Cells.Item(35, 5).Value = 7.48
$wsReco.Cells.Item(35, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(35, 7).Value = '👀 À surveiller'

# Row 36: TRACTAFRIC MOTORS CI (PRSC)
$wsReco.Cells.Item(36, 1).Value = 'TRACTAFRIC MOTORS CI (PRSC)'
$wsReco.Cells.Item(36, 2).Value = 1
$wsReco.Cells.Item(36, 3).Value = 1
$wsReco.Cells.Item(36, 4).Value = -0.02
$wsReco.Cells.Item(36, 5).Value = 4.16
$wsReco.Cells.Item(36, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(36, 7).Value = '👀 À surveiller'

# Row 37: AFRICA GLOBAL LOGISTICS CI (SDSC)
$wsReco.Cells.Item(37, 1).Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$wsReco.Cells.Item(37, 2).Value = 1
$wsReco.Cells.Item(37, 3).Value = 1
$wsReco.Cells.Item(37, 4).Value = -0.55
$wsReco.Cells.Item(37, 5).Value = -3.99
$wsReco.Cells.Item(37, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(37, 7).Value = '👀 À surveiller'

# Row 38: TOTALENERGIES MARKETING SN (TTLS)
$wsReco.Cells.Item(38, 1).Value = 'TOTALENERGIES MARKETING SN (TTLS)'
$wsReco.Cells.Item(38, 2).Value = 1
$wsReco.Cells.Item(38, 3).Value = 1
$wsReco.Cells.Item(38, 4).Value = -0.62
$wsReco.Cells.Item(38, 5).Value = 3.02
$wsReco.Cells.Item(38, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(38, 7).Value = '👀 À surveiller'

# Row 39: ORAGROUP TOGO (ORGT)
$wsReco.Cells.Item(39, 1).Value = 'ORAGROUP TOGO (ORGT)'
$wsReco.Cells.Item(39, 2).Value = 0
$wsReco.Cells.Item(39, 3).Value = 1
$wsReco.Cells.Item(39, 4).Value = -2.42
$wsReco.Cells.Item(39, 5).Value = -2.42
$wsReco.Cells.Item(39, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(39, 7).Value = '➖ Neutre'

# Row 40: SMB CI (SMBC)
$wsReco.Cells.Item(40, 1).Value = 'SMB CI (SMBC)'
$wsReco.Cells.Item(40, 2).Value = 1
$wsReco.Cells.Item(40, 3).Value = 1
$wsReco.Cells.Item(40, 4).Value = -2.47
$wsReco.Cells.Item(40, 5).Value = -6.21
$wsReco.Cells.Item(40, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(40, 7).Value = '👀 À surveiller'

# Row 41: NEI-CEDA CI (NEIC)
$wsReco.Cells.Item(41, 1).Value = 'NEI-CEDA CI (NEIC)'
$wsReco.Cells.Item(41, 2).Value = 0
$wsReco.Cells.Item(41, 3).Value = 1
$wsReco.Cells.Item(41, 4).Value = -2.5
$wsReco.Cells.Item(41, 5).Value = -2.5
$wsReco.Cells.Item(41, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(41, 7).Value = '➖ Neutre'

# Row 42: ORANGE COTE D'IVOIRE (ORAC)
$wsReco.Cells.Item(42, 1).Value = 'ORANGE COTE D''IVOIRE (ORAC)'
$wsReco.Cells.Item(42, 2).Value = 0
$wsReco.Cells.Item(42, 3).Value = 1
$wsReco.Cells.Item(42, 4).Value = -3.16
$wsReco.Cells.Item(42, 5).Value = -3.16
$wsReco.Cells.Item(42, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(42, 7).Value = '➖ Neutre'

# Row 43: BANK OF AFRICA BF (BOABF)
$wsReco.Cells.Item(43, 1).Value = 'BANK OF AFRICA BF (BOABF)'
$wsReco.Cells.Item(43, 2).Value = 0
$wsReco.Cells.Item(43, 3).Value = 1
$wsReco.Cells.Item(43, 4).Value = -5.33
$wsReco.Cells.Item(43, 5).Value = -5.33
$wsReco.Cells.Item(43, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(43, 7).Value = '➖ Neutre'

# Row 44: SODE CI (SDCC)
$wsReco.Cells.Item(44, 1).Value = 'SODE CI (SDCC)'
$wsReco.Cells.Item(44, 2).Value = 0
$wsReco.Cells.Item(44, 3).Value = 1
$wsReco.Cells.Item(44, 4).Value = -5.4
$wsReco.Cells.Item(44, 5).Value = -5.4
$wsReco.Cells.Item(44, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(44, 7).Value = '➖ Neutre'

# Row 45: VIVO ENERGY CI (SHEC)
$wsReco.Cells.Item(45, 1).Value = 'VIVO ENERGY CI (SHEC)'
$wsReco.Cells.Item(45, 2).Value = 0
$wsReco.Cells.Item(45, 3).Value = 2
$wsReco.Cells.Item(45, 4).Value = -7
$wsReco.Cells.Item(45, 5).Value = -3.79
$wsReco.Cells.Item(45, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(45, 7).Value = '➖ Neutre'

# Row 46: FILTISAC CI (FTSC)
$wsReco.Cells.Item(46, 1).Value = 'FILTISAC CI (FTSC)'
$wsReco.Cells.Item(46, 2).Value = 0
$wsReco.Cells.Item(46, 3).Value = 2
$wsReco.Cells.Item(46, 4).Value = -9.210000000000001
$wsReco.Cells.Item(46, 5).Value = -5.66
$wsReco.Cells.Item(46, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(46, 7).Value = '➖ Neutre'

# Row 47: SERVAIR ABIDJAN CI (ABJC)
$wsReco.Cells.Item(47, 1).Value = 'SERVAIR ABIDJAN CI (ABJC)'
$wsReco.Cells.Item(47, 2).Value = 0
$wsReco.Cells.Item(47, 3).Value = 2
$wsReco.Cells.Item(47, 4).Value = -14.92
$wsReco.Cells.Item(47, 5).Value = -7.45
$wsReco.Cells.Item(47, 6).Value = '🟡 Observer'
$wsReco.Cells.Item(47, 7).Value = '➖ Neutre'

# Remove the old trailing row (sheet now spans A1:G47)
$wsReco.Rows(48).Delete()

$wsYtd = $wb.Worksheets.Item("Top_YTD")

# Row 2: BRVM - SERVICES PUBLICS
$wsYtd.Cells.Item(2, 1).Value = 'BRVM - SERVICES PUBLICS'
$wsYtd.Cells.Item(2, 2).Value = 10557915.2

# Row 3: SAFCA CI
$wsYtd.Cells.Item(3, 1).Value = 'SAFCA CI'
$wsYtd.Cells.Item(3, 2).Value = 382036.43

# Row 4: CFAO MOTORS CI
$wsYtd.Cells.Item(4, 1).Value = 'CFAO MOTORS CI'
$wsYtd.Cells.Item(4, 2).Value = 362932.28

# Row 5: BRVM - AUTRES SECTEURS
$wsYtd.Cells.Item(5, 1).Value = 'BRVM - AUTRES SECTEURS'
$wsYtd.Cells.Item(5, 2).Value = 341049.57

# Row 6: UNIWAX CI
$wsYtd.Cells.Item(6, 1).Value = 'UNIWAX CI'
$wsYtd.Cells.Item(6, 2).Value = 229856.3

# Row 7: NEI-CEDA CI
$wsYtd.Cells.Item(7, 1).Value = 'NEI-CEDA CI'
$wsYtd.Cells.Item(7, 2).Value = 228177.96

# Row 8: AIR LIQUIDE CI
$wsYtd.Cells.Item(8, 1).Value = 'AIR LIQUIDE CI'
$wsYtd.Cells.Item(8, 2).Value = 170303.84

# Row 9: SETAO CI
$wsYtd.Cells.Item(9, 1).Value = 'SETAO CI'
$wsYtd.Cells.Item(9, 2).Value = 163699.04

# Row 10: BRVM - TRANSPORT
$wsYtd.Cells.Item(10, 1).Value = 'BRVM - TRANSPORT'
$wsYtd.Cells.Item(10, 2).Value = 57790.3

# Row 11: BRVM - DISTRIBUTION
$wsYtd.Cells.Item(11, 1).Value = 'BRVM - DISTRIBUTION'
$wsYtd.Cells.Item(11, 2).Value = 48895.09

Write-Host "BRVM recommandations + Top_YTD refreshed"
